$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.658.02"
$ws.Range("E2").Value = "  +3.07%  "
$ws.Range("D3").Value = "2.424.69"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.47"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +12.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.62"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.648"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.677"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +12.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.56"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0950"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.68"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("E13").Value = "  +3.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.43"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +17.31%  "
$ws.Range("E15").Value = "  +3.02%  "
$ws.Range("D16").Value = "2.790.12"
$ws.Range("E16").Value = "  +9.02%  "
$ws.Range("D17").Value = "2.428.15"
$ws.Range("E17").Value = "  +9.34%  "
$ws.Range("D18").Value = "43.674.39"
$ws.Range("E18").Value = "  +3.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000111"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.52"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +5.59%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.40"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.09%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "275.67"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +19.99%  "
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.52"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.47"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.65"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +7.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.10"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +6.12%  "
$ws.Range("E27").Value = "  -0.11%  "
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "23.16"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +11.18%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "178.82"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.14%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.23"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "38.30"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.30%  "
$ws.Range("E33").Value = "  +4.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0942"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +7.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.03"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +7.70%  "
$ws.Range("E36").Value = "  +6.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.92"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0373"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.05"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.109"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.52%  "
$ws.Range("E41").Value = "  +21.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.64"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +26.02%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "129.88"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +28.59%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.237"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "70.38"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.55%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.78"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.30%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.78"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +16.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.74"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +7.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "87.13"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +69.00%  "
$ws.Range("E51").Value = "  +4.90%  "
